{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the practice\n// table with the newly generated set of problems/answers.\n// Each old string is unique within the document, so a direct search +\n// replace (first/only match) is safe and preserves the existing run\n// formatting (font, size) of each cell.\nconst replacements = [\n  [\"30\u00f76=5, 0\", \"17\u00f77=2, 3\"],\n  [\"10\u00f76=1, 4\", \"11\u00f74=2, 3\"],\n  [\"82\u00f74=20, 2\", \"42\u00f72=21, 0\"],\n  [\"80\u00f76=13, 2\", \"64\u00f79=7, 1\"],\n  [\"31\u00f77=4, 3\", \"43\u00f75=8, 3\"],\n  [\"24\u00f76=4, 0\", \"28\u00f76=4, 4\"],\n  [\"79\u00f72=39, 1\", \"14\u00f75=2, 4\"],\n  [\"73\u00f75=14, 3\", \"99\u00f79=11, 0\"],\n  [\"77\u00f74=19, 1\", \"29\u00f73=9, 2\"],\n  [\"80\u00f75=16, 0\", \"63\u00f73=21, 0\"],\n  [\"11\u00f72=5, 1\", \"24\u00f75=4, 4\"],\n  [\"61\u00f75=12, 1\", \"90\u00f75=18, 0\"],\n  [\"62\u00f76=10, 2\", \"56\u00f72=28, 0\"],\n  [\"56\u00f79=6, 2\", \"51\u00f79=5, 6\"],\n  [\"30\u00f75=6, 0\", \"60\u00f74=15, 0\"],\n  [\"91\u00f75=18, 1\", \"44\u00f76=7, 2\"],\n  [\"39\u00f72=19, 1\", \"94\u00f72=47, 0\"],\n  [\"10\u00f79=1, 1\", \"55\u00f72=27, 1\"],\n  [\"55\u00f74=13, 3\", \"63\u00f75=12, 3\"],\n  [\"75\u00f74=18, 3\", \"51\u00f79=5, 6\"],\n  [\"41\u00f75=8, 1\", \"92\u00f74=23, 0\"],\n  [\"12\u00f78=1, 4\", \"55\u00f79=6, 1\"],\n  [\"44\u00f73=14, 2\", \"22\u00f77=3, 1\"],\n  [\"76\u00f77=10, 6\", \"45\u00f78=5, 5\"],\n  [\"35\u00f77=5, 0\", \"21\u00f73=7, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the practice\n# table with the newly generated set of problems/answers. Each old string\n# is unique within the document, so Find/Replace (ReplaceAll) targets the\n# correct cell every time and preserves the existing run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"30\u00f76=5, 0\", \"17\u00f77=2, 3\"),\n  @(\"10\u00f76=1, 4\", \"11\u00f74=2, 3\"),\n  @(\"82\u00f74=20, 2\", \"42\u00f72=21, 0\"),\n  @(\"80\u00f76=13, 2\", \"64\u00f79=7, 1\"),\n  @(\"31\u00f77=4, 3\", \"43\u00f75=8, 3\"),\n  @(\"24\u00f76=4, 0\", \"28\u00f76=4, 4\"),\n  @(\"79\u00f72=39, 1\", \"14\u00f75=2, 4\"),\n  @(\"73\u00f75=14, 3\", \"99\u00f79=11, 0\"),\n  @(\"77\u00f74=19, 1\", \"29\u00f73=9, 2\"),\n  @(\"80\u00f75=16, 0\", \"63\u00f73=21, 0\"),\n  @(\"11\u00f72=5, 1\", \"24\u00f75=4, 4\"),\n  @(\"61\u00f75=12, 1\", \"90\u00f75=18, 0\"),\n  @(\"62\u00f76=10, 2\", \"56\u00f72=28, 0\"),\n  @(\"56\u00f79=6, 2\", \"51\u00f79=5, 6\"),\n  @(\"30\u00f75=6, 0\", \"60\u00f74=15, 0\"),\n  @(\"91\u00f75=18, 1\", \"44\u00f76=7, 2\"),\n  @(\"39\u00f72=19, 1\", \"94\u00f72=47, 0\"),\n  @(\"10\u00f79=1, 1\", \"55\u00f72=27, 1\"),\n  @(\"55\u00f74=13, 3\", \"63\u00f75=12, 3\"),\n  @(\"75\u00f74=18, 3\", \"51\u00f79=5, 6\"),\n  @(\"41\u00f75=8, 1\", \"92\u00f74=23, 0\"),\n  @(\"12\u00f78=1, 4\", \"55\u00f79=6, 1\"),\n  @(\"44\u00f73=14, 2\", \"22\u00f77=3, 1\"),\n  @(\"76\u00f77=10, 6\", \"45\u00f78=5, 5\"),\n  @(\"35\u00f77=5, 0\", \"21\u00f73=7, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
